$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Replace the stray "(2012, Mandić, Vladimir)" citation token with
#        "mandic2012a" inside the big citation-list cells that still contain it.
$citeCells = @("J41", "L11", "L40", "L44", "L53")
foreach ($addr in $citeCells) {
    $cell = $ws.Range($addr)
    $old = $cell.Text
    $new = $old.Replace("(2012, Mandić, Vladimir)", "mandic2012a")
    $cell.Value = $new
}

# --- 2) Update the three remaining column-U cells whose text drops the
#        "becker2008a"/"trienekens2009a" token combination.
$ws.Range("U5").Value = "trienekens2009a"
$ws.Range("U9").Value = "becker2008a, birkhölzer2011a"
$ws.Range("U49").Value = "trienekens2009a"

# --- 3) Clear out column U for every other data row (3-53, excluding the
#        header row 2 and rows 5, 9, 49 handled above).
$rowsToClear = @(3,4,6,7,8,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,50,51,52,53)
foreach ($r in $rowsToClear) {
    $ws.Range("U$r").ClearContents()
}
